# bug fix in 191
# Adds newly-recorded sprint rows to the AMSIN and AMS history sheets, and
# fixes up the AMS row for 189livewp (style normalisation + a corrected run
# timestamp) as part of the 191 bug fix.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "AMSIN": append rows 45-50
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

$amsinRows = @(
    @{ Row=45; A="2024-03-28"; B=45379.53211018519;  C="190fstwp";    D=46; E=46; F=0; G=0.68 },
    @{ Row=46; A="2024-03-29"; B=45380.4774421875;   C="190scndwp";   D=46; E=46; F=0; G=0.71 },
    @{ Row=47; A="2024-04-01"; B=45383.36922920139;  C="190fnlwp";    D=46; E=46; F=0; G=0.67 },
    @{ Row=48; A="2024-05-02"; B=45414.47397732639;  C="191fstwp";    D=46; E=41; F=5; G=1.89 },
    @{ Row=49; A="2024-05-02"; B=45414.55501819444;  C="191wprofile"; D=46; E=46; F=0; G=0.87 },
    @{ Row=50; A="2024-05-03"; B=45415.33852509259;  C="191lstwp";    D=46; E=46; F=0; G=0.73 }
)

foreach ($r in $amsinRows) {
    $row = $r.Row
    # Column A holds a date-like string; a leading apostrophe keeps it text
    # (matching the existing rows above it) instead of auto-converting to a
    # date serial.
    $wsAmsin.Range("A$row").Value = "'" + $r.A
    $wsAmsin.Range("B$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsAmsin.Range("B$row").Value = $r.B
    $wsAmsin.Range("C$row").Value = $r.C
    $wsAmsin.Range("D$row").Value = $r.D
    $wsAmsin.Range("E$row").Value = $r.E
    $wsAmsin.Range("F$row").Value = $r.F
    $wsAmsin.Range("G$row").Value = $r.G
}

# ---------------------------------------------------------------------
# Sheet "AMS": fix row 36 (189livewp) and append rows 37-39
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Normalise row 36's formatting to match row 35 (it previously had no
# explicit style while its neighbours did) and correct its run timestamp.
$wsAms.Range("A35:G35").Copy()
$wsAms.Range("A36:G36").PasteSpecial(-4122)
$wsAms.Application.CutCopyMode = $false
$wsAms.Range("B36").Value = 45359.73304665509

$amsRows = @(
    @{ Row=37; A="2024-04-01"; B=45383.53188237268;  C="190betawpp"; D=46; E=46; F=0; G=0.86 },
    @{ Row=38; A="2024-04-01"; B=45383.85208363426;  C="190livewp";  D=46; E=44; F=2; G=1.23 },
    @{ Row=39; A="2024-05-03"; B=45415.60713529508;  C="191betyawp"; D=46; E=46; F=0; G=0.82 }
)

foreach ($r in $amsRows) {
    $row = $r.Row
    $wsAms.Range("A$row").Value = "'" + $r.A
    $wsAms.Range("B$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsAms.Range("B$row").Value = $r.B
    $wsAms.Range("C$row").Value = $r.C
    $wsAms.Range("D$row").Value = $r.D
    $wsAms.Range("E$row").Value = $r.E
    $wsAms.Range("F$row").Value = $r.F
    $wsAms.Range("G$row").Value = $r.G
}

Write-Host "edit applied"
